# "Created logic for building the alignments"
#
# - Rename the worksheet tab from "sequences.txt" to "lassa-data".
# - Fill in the remaining columns (B:H) for the last data row (1937),
#   which previously only had column A populated.
# - Move the active-cell selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab.
$ws.Name = "lassa-data"

# Row 1937 only had "Accession" (column A) filled in; fill the rest of the
# row with the usual placeholder "-" values, and mark patent_related (G) true.
$ws.Range("B1937").Value = "-"
$ws.Range("C1937").Value = "-"
$ws.Range("D1937").Value = "-"

$ws.Range("E1937").Value = "-"
# Column E normally carries the sheet's left-aligned column style; this new
# row's cell should stay on the default/Normal style.
$ws.Range("E1937").Style = "Normal"

$ws.Range("F1937").Value = "-"
$ws.Range("G1937").Value = $true
$ws.Range("H1937").Value = "-"

# Move the selection to A2.
$ws.Range("A2").Select() | Out-Null
